$codonData = @(
    @{Row=1; B=0; C="A"; D="Gly"},
    @{Row=2; B=1; C="A"; D="Gly"},
    @{Row=3; B=2; C="A"; D="Gly"},
    @{Row=4; B=3; C="A"; D="Gly"},
    @{Row=5; B=4; C="B"; D="Glu"},
    @{Row=6; B=5; C="B"; D="Glu"},
    @{Row=7; B=6; C="C"; D="Asp"},
    @{Row=8; B=7; C="C"; D="Asp"},
    @{Row=9; B=8; C="D"; D="Val"},
    @{Row=10; B=9; C="D"; D="Val"},
    @{Row=11; B=10; C="D"; D="Val"},
    @{Row=12; B=11; C="D"; D="Val"},
    @{Row=13; B=12; C="E"; D="Ala"},
    @{Row=14; B=13; C="E"; D="Ala"},
    @{Row=15; B=14; C="E"; D="Ala"},
    @{Row=16; B=15; C="E"; D="Ala"},
    @{Row=17; B=16; C="F"; D="Arg"},
    @{Row=18; B=17; C="F"; D="Arg"},
    @{Row=19; B=18; C="G"; D="Ser"},
    @{Row=20; B=19; C="G"; D="Ser"},
    @{Row=21; B=20; C="H"; D="Lys"},
    @{Row=22; B=21; C="H"; D="Lys"},
    @{Row=23; B=22; C="I"; D="Asn"},
    @{Row=24; B=23; C="I"; D="Asn"},
    @{Row=25; B=24; C="J"; D="Met"},
    @{Row=26; B=25; C="K"; D="Ile"},
    @{Row=27; B=26; C="K"; D="Ile"},
    @{Row=28; B=27; C="K"; D="Ile"},
    @{Row=29; B=28; C="L"; D="Thr"},
    @{Row=30; B=29; C="L"; D="Thr"},
    @{Row=31; B=30; C="L"; D="Thr"},
    @{Row=32; B=31; C="L"; D="Thr"},
    @{Row=33; B=32; C="M"; D="Trp"},
    @{Row=34; B=33; C="N"; D="Stop"},
    @{Row=35; B=34; C="O"; D="Cys"},
    @{Row=36; B=35; C="O"; D="Cys"},
    @{Row=37; B=36; C="P"; D="Tyr"},
    @{Row=38; B=37; C="P"; D="Tyr"},
    @{Row=39; B=38; C="Q"; D="Leu"},
    @{Row=40; B=39; C="Q"; D="Leu"},
    @{Row=41; B=40; C="R"; D="Phe"},
    @{Row=42; B=41; C="R"; D="Phe"},
    @{Row=43; B=42; C="G"; D="Ser"},
    @{Row=44; B=43; C="G"; D="Ser"},
    @{Row=45; B=44; C="G"; D="Ser"},
    @{Row=46; B=45; C="G"; D="Ser"},
    @{Row=47; B=46; C="F"; D="Arg"},
    @{Row=48; B=47; C="F"; D="Arg"},
    @{Row=49; B=48; C="F"; D="Arg"},
    @{Row=50; B=49; C="F"; D="Arg"},
    @{Row=51; B=50; C="S"; D="Gln"},
    @{Row=52; B=51; C="S"; D="Gln"},
    @{Row=53; B=52; C="T"; D="His"},
    @{Row=54; B=53; C="T"; D="His"},
    @{Row=55; B=54; C="Q"; D="Leu"},
    @{Row=56; B=55; C="Q"; D="Leu"},
    @{Row=57; B=56; C="Q"; D="Leu"},
    @{Row=58; B=57; C="Q"; D="Leu"},
    @{Row=59; B=58; C="U"; D="Pro"},
    @{Row=60; B=59; C="U"; D="Pro"},
    @{Row=61; B=60; C="U"; D="Pro"},
    @{Row=62; B=61; C="U"; D="Pro"}
)


$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($item in $codonData) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B        # Column B
    $ws.Cells.Item($r, 3).Value = $item.C        # Column C (one-letter amino acid code)
    $ws.Cells.Item($r, 4).Value = $item.D        # Column D (three-letter amino acid code)
}

# Update the G44 formula to include E19 and E20 in the sum (synonymous codon optimization)
$ws.Range("G44").Formula = "=E45+E46+E44+E43+E19+E20"

# Update the sheet view to match saved state (scrolled down, different selection)
$ws.Range("C63").Select()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
